$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 30 (shifts MS-M and everything below down by one)
$ws.Rows.Item(30).Insert()

# Populate the new row 30 with the MS-F data
$ws.Range("A30").Value = "MS-F"
$ws.Range("B30").Value = -4.86
$ws.Range("C30").Value = -3.442
$ws.Range("D30").Value = 13.288
$ws.Range("E30").Value = -2665.761431
$ws.Range("F30").Value = -1.129
$ws.Range("G30").Value = 727.90863
$ws.Range("H30").Value = 831.593118
$ws.Range("I30").Value = 671
$ws.Range("J30").Value = 0.321307
